$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.972.73"
$ws.Range("E2").Value = "  +0.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.311.76"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "541.92"
$ws.Range("E5").Value = "  -0.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.99"
$ws.Range("E6").Value = "  -1.72%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("E8").Value = "  -2.14%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.309.31"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  -0.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.54"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.32"
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.722.84"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.960.13"
$ws.Range("E16").Value = "  +1.22%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000132"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.301.33"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.47"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.07"
$ws.Range("E20").Value = "  -2.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "312.15"
$ws.Range("E21").Value = "  -0.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.54"
$ws.Range("E22").Value = "  -0.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.30"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  -2.94%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.75"
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("E28").Value = "  +2.82%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.56"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.18"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("E31").Value = "  -1.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0724"
$ws.Range("E32").Value = "  -2.12%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.84"
$ws.Range("E33").Value = "  -0.69%  "
$ws.Range("E34").Value = "  +0.87%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.378"
$ws.Range("E35").Value = "  -2.07%  "
$ws.Range("E36").Value = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.69"
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.99"
$ws.Range("E39").Value = "  -2.08%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "316.65"
$ws.Range("E40").Value = "  +0.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "37.70"
$ws.Range("E41").Value = "  -0.75%  "
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "135.93"
$ws.Range("E43").Value = "  -3.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.43"
$ws.Range("E44").Value = "  -0.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0941"
$ws.Range("E45").Value = "  -1.37%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.569"
$ws.Range("E46").Value = "  +1.72%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "18.71"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0489"
$ws.Range("E48").Value = "  -1.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₆0223"
$ws.Range("E49").Value = "  +21.23%  "
$ws.Range("E50").Value = "  +0.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.00"
$ws.Range("E51").Value = "  -0.07%  "
